$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add a new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the look of the header row used on the other sheets (bold, centered,
# thin border) and the date format used in column A.
$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$wsForecast.Range("A2:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Data rows
$data = @(
    @(45326.99999999999, 63, 54.54973733491261, 71.37388121536011),
    @(45333.99999999999, 51, 43.03209703737401, 59.45482045804636),
    @(45340.99999999999, 39, 30.17736739073011, 47.98616072417735),
    @(45347.99999999999, 27, 17.58568250501225, 35.50499647550254),
    @(45354.99999999999, 15, 6.951903911116219, 23.98136902536054),
    @(45361.99999999999, 3, -5.49979768087668, 12.20190063087426),
    @(45368.99999999999, 0, -17.63425925550419, -0.1212872122308769),
    @(45375.99999999999, 0, -28.60952008805263, -12.71266015761912),
    @(45382.99999999999, 0, -41.20755522114779, -24.46689010123373),
    @(45389.99999999999, 0, -53.21005483726067, -36.15431054013805),
    @(45396.99999999999, 0, -65.30886919857355, -48.55701286813434),
    @(45403.99999999999, 0, -77.36133281376085, -60.82063637643504)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

Write-Host "PO Forecast sheet added and headers renamed."
